$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "RelivePos" text values from "186,6.89,88" to "186,0,88"
# in every cell that currently holds that value (E2, E3, E4).
$ws.Range("E2").Value = "186,0,88"
$ws.Range("E3").Value = "186,0,88"
$ws.Range("E4").Value = "186,0,88"

# Update the current selection to F7, matching the saved sheet view state.
$ws.Range("F7").Select()
